{"js": "// Replace the first paragraph's split \"!@\" / \"test@!\" runs (with the\n// gramStart/gramEnd proofing-error markers between them) with a single\n// merged run \"!@test@!\", then add a blank paragraph followed by a new\n// paragraph containing \"!@choice@!\" (multi-select choice support).\nconst body = context.document.body;\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the \"!@test@!\" paragraph defensively (it is paragraph 0 in this\n// document) instead of assuming a fixed index.\nlet firstParagraph = body.paragraphs.items.find(p => p.text === \"!@test@!\");\nif (!firstParagraph) firstParagraph = body.paragraphs.items[0];\n\n// Build the replacement content as a flat-OPC OOXML package fragment so we\n// can fully control the resulting markup: one merged run in paragraph 1\n// (dropping the stray <w:proofErr/> markers), a truly empty paragraph, and\n// a new paragraph carrying the \"!@choice@!\" text.\nconst ns = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\nconst newBodyXml =\n  '<w:p><w:r><w:t>!@test@!</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t>!@choice@!</w:t></w:r></w:p>';\nconst ooxmlPackage =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document ' + ns + '><w:body>' + newBodyXml + '</w:body></w:document></pkg:xmlData>' +\n  '</pkg:part></pkg:package>';\n\nfirstParagraph.insertOoxml(ooxmlPackage, \"Replace\");\nawait context.sync();\n", "ps1": "# Continued with implementing multi select support.\n# Merge the first paragraph's split \"!@\" / \"test@!\" runs (with the stray\n# gramStart/gramEnd proofing-error markers between them) into a single run\n# reading \"!@test@!\", then append a blank paragraph followed by a new\n# paragraph containing \"!@choice@!\".\n\n$d = $word.ActiveDocument\n\n# Clear the first paragraph's content (text + proofErr markers) but keep\n# the paragraph mark itself, then retype the merged text into it so it\n# collapses back down to a single run.\n$firstPara = $d.Paragraphs(1)\n$firstPara.Range.Delete()\n$firstPara.Range.InsertAfter(\"!@test@!\")\n\n# Append the new blank paragraph + the \"!@choice@!\" paragraph right after\n# the (now merged) first paragraph, in one shot via InsertXML so the blank\n# paragraph ends up with no leftover empty run. The target position is\n# just before the first paragraph's own paragraph mark.\n$insertAt = $firstPara.Range.End - 1\n$endRange = $d.Range($insertAt, $insertAt)\n$endRange.InsertXML('<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p/><w:p><w:r><w:t>!@choice@!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')\n"}
